# Adds rows 32-38 to the "Artfynd" sheet (species-observation records),
# matching the upstream commit that appended 7 new fynd (finds) rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Range("A32").Value = 112203732
$ws.Range("B32").Value = 90658
$ws.Range("C32").Value = "Ovaliderad"
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 4361
$ws.Range("F32").Value = "Orange taggsvamp"
$ws.Range("G32").Value = "Hydnellum aurantiacum"
$ws.Range("H32").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I32").Value = "'30"
$ws.Range("I32").ClearFormats()
$ws.Range("J32").Value = "fruktkroppar"
$ws.Range("K32").Value = "'"
$ws.Range("K32").ClearFormats()
$ws.Range("N32").Value = "'"
$ws.Range("N32").ClearFormats()
$ws.Range("P32").Value = "A 30779, Storön, Sm"
$ws.Range("Q32").Value = 594803.3283005389
$ws.Range("R32").Value = 6396140.582550677
$ws.Range("S32").Value = 10
$ws.Range("T32").Value = "Kalmar"
$ws.Range("U32").Value = "Västervik"
$ws.Range("V32").Value = "Småland"
$ws.Range("W32").Value = "Gladhammar"
$ws.Range("Y32").Value = "'2023-09-19"
$ws.Range("Y32").ClearFormats()
$ws.Range("Z32").Value = "00:00"
$ws.Range("AA32").Value = "'2023-09-19"
$ws.Range("AA32").ClearFormats()
$ws.Range("AB32").Value = "00:00"
$ws.Range("AD32").Value = $false
$ws.Range("AE32").Value = $false
$ws.Range("AF32").Value = "'"
$ws.Range("AF32").ClearFormats()
$ws.Range("AG32").Value = $false
$ws.Range("AT32").Value = "'"
$ws.Range("AT32").ClearFormats()
$ws.Range("AW32").Value = "Magnus Kasselstrand"
$ws.Range("AX32").Value = "Magnus Kasselstrand, Ingvor Kasselstrand, Gunilla Nilsson, Larsgunnar Nilsson"
$ws.Range("AY32").Value = "'"
$ws.Range("AY32").ClearFormats()

# Row 33
$ws.Range("A33").Value = 112204167
$ws.Range("B33").Value = 90689
$ws.Range("C33").Value = "Ovaliderad"
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 5966
$ws.Range("F33").Value = "Motaggsvamp"
$ws.Range("G33").Value = "Sarcodon squamosus"
$ws.Range("H33").Value = "(Schaeff.) Quél."
$ws.Range("I33").Value = "'10"
$ws.Range("I33").ClearFormats()
$ws.Range("J33").Value = "fruktkroppar"
$ws.Range("K33").Value = "'"
$ws.Range("K33").ClearFormats()
$ws.Range("N33").Value = "'"
$ws.Range("N33").ClearFormats()
$ws.Range("P33").Value = "A 30799, Storön, Sm"
$ws.Range("Q33").Value = 594925.0145863529
$ws.Range("R33").Value = 6396227.912029894
$ws.Range("S33").Value = 100
$ws.Range("T33").Value = "Kalmar"
$ws.Range("U33").Value = "Västervik"
$ws.Range("V33").Value = "Småland"
$ws.Range("W33").Value = "Gladhammar"
$ws.Range("Y33").Value = "'2023-09-19"
$ws.Range("Y33").ClearFormats()
$ws.Range("Z33").Value = "00:00"
$ws.Range("AA33").Value = "'2023-09-19"
$ws.Range("AA33").ClearFormats()
$ws.Range("AB33").Value = "00:00"
$ws.Range("AC33").Value = "2+7+1"
$ws.Range("AD33").Value = $false
$ws.Range("AE33").Value = $false
$ws.Range("AF33").Value = "'"
$ws.Range("AF33").ClearFormats()
$ws.Range("AG33").Value = $false
$ws.Range("AT33").Value = "'"
$ws.Range("AT33").ClearFormats()
$ws.Range("AW33").Value = "Magnus Kasselstrand"
$ws.Range("AX33").Value = "Magnus Kasselstrand, Ingvor Kasselstrand, Gunilla Nilsson, Larsgunnar Nilsson"
$ws.Range("AY33").Value = "'"
$ws.Range("AY33").ClearFormats()

# Row 34
$ws.Range("A34").Value = 112204311
$ws.Range("B34").Value = 96348
$ws.Range("C34").Value = "Ovaliderad"
$ws.Range("D34").Value = "VU"
$ws.Range("E34").Value = 220787
$ws.Range("F34").Value = "Knärot"
$ws.Range("G34").Value = "Goodyera repens"
$ws.Range("H34").Value = "(L.) R. Br."
$ws.Range("I34").Value = "'300"
$ws.Range("I34").ClearFormats()
$ws.Range("J34").Value = "plantor/tuvor"
$ws.Range("K34").Value = "'"
$ws.Range("K34").ClearFormats()
$ws.Range("L34").Value = "'"
$ws.Range("L34").ClearFormats()
$ws.Range("N34").Value = "'"
$ws.Range("N34").ClearFormats()
$ws.Range("P34").Value = "Storön, Samsvik, Sm"
$ws.Range("Q34").Value = 595112.440808123
$ws.Range("R34").Value = 6396025.415650645
$ws.Range("S34").Value = 5
$ws.Range("T34").Value = "Kalmar"
$ws.Range("U34").Value = "Västervik"
$ws.Range("V34").Value = "Småland"
$ws.Range("W34").Value = "Gladhammar"
$ws.Range("Y34").Value = "'2023-09-19"
$ws.Range("Y34").ClearFormats()
$ws.Range("Z34").Value = "00:00"
$ws.Range("AA34").Value = "'2023-09-19"
$ws.Range("AA34").ClearFormats()
$ws.Range("AB34").Value = "00:00"
$ws.Range("AD34").Value = $false
$ws.Range("AE34").Value = $false
$ws.Range("AF34").Value = "'"
$ws.Range("AF34").ClearFormats()
$ws.Range("AG34").Value = $false
$ws.Range("AT34").Value = "'"
$ws.Range("AT34").ClearFormats()
$ws.Range("AW34").Value = "Larsgunnar Nilsson"
$ws.Range("AX34").Value = "Larsgunnar Nilsson, Gunilla Nilsson, Ingvor Kasselstrand, Magnus Kasselstrand"
$ws.Range("AY34").Value = "'"
$ws.Range("AY34").ClearFormats()

# Row 35
$ws.Range("A35").Value = 112203737
$ws.Range("B35").Value = 90709
$ws.Range("C35").Value = "Ovaliderad"
$ws.Range("D35").Value = "NT"
$ws.Range("E35").Value = 5448
$ws.Range("F35").Value = "Svartvit taggsvamp"
$ws.Range("G35").Value = "Phellodon connatus"
$ws.Range("H35").Value = "(Schultz) nom.prov"
$ws.Range("I35").Value = "'10"
$ws.Range("I35").ClearFormats()
$ws.Range("J35").Value = "fruktkroppar"
$ws.Range("K35").Value = "'"
$ws.Range("K35").ClearFormats()
$ws.Range("N35").Value = "'"
$ws.Range("N35").ClearFormats()
$ws.Range("P35").Value = "A 30779, Storön, Sm"
$ws.Range("Q35").Value = 594803.3283005389
$ws.Range("R35").Value = 6396140.582550677
$ws.Range("S35").Value = 10
$ws.Range("T35").Value = "Kalmar"
$ws.Range("U35").Value = "Västervik"
$ws.Range("V35").Value = "Småland"
$ws.Range("W35").Value = "Gladhammar"
$ws.Range("Y35").Value = "'2023-09-19"
$ws.Range("Y35").ClearFormats()
$ws.Range("Z35").Value = "00:00"
$ws.Range("AA35").Value = "'2023-09-19"
$ws.Range("AA35").ClearFormats()
$ws.Range("AB35").Value = "00:00"
$ws.Range("AD35").Value = $false
$ws.Range("AE35").Value = $false
$ws.Range("AF35").Value = "'"
$ws.Range("AF35").ClearFormats()
$ws.Range("AG35").Value = $false
$ws.Range("AT35").Value = "'"
$ws.Range("AT35").ClearFormats()
$ws.Range("AW35").Value = "Magnus Kasselstrand"
$ws.Range("AX35").Value = "Magnus Kasselstrand, Ingvor Kasselstrand, Gunilla Nilsson, Larsgunnar Nilsson"
$ws.Range("AY35").Value = "'"
$ws.Range("AY35").ClearFormats()

# Row 36
$ws.Range("A36").Value = 112204297
$ws.Range("B36").Value = 96348
$ws.Range("C36").Value = "Ovaliderad"
$ws.Range("D36").Value = "VU"
$ws.Range("E36").Value = 220787
$ws.Range("F36").Value = "Knärot"
$ws.Range("G36").Value = "Goodyera repens"
$ws.Range("H36").Value = "(L.) R. Br."
$ws.Range("I36").Value = "'10"
$ws.Range("I36").ClearFormats()
$ws.Range("J36").Value = "plantor/tuvor"
$ws.Range("K36").Value = "'"
$ws.Range("K36").ClearFormats()
$ws.Range("L36").Value = "'"
$ws.Range("L36").ClearFormats()
$ws.Range("N36").Value = "'"
$ws.Range("N36").ClearFormats()
$ws.Range("P36").Value = "Storön, Samsvik, Sm"
$ws.Range("Q36").Value = 595097.205446711
$ws.Range("R36").Value = 6396058.205299424
$ws.Range("S36").Value = 5
$ws.Range("T36").Value = "Kalmar"
$ws.Range("U36").Value = "Västervik"
$ws.Range("V36").Value = "Småland"
$ws.Range("W36").Value = "Gladhammar"
$ws.Range("Y36").Value = "'2023-09-19"
$ws.Range("Y36").ClearFormats()
$ws.Range("Z36").Value = "00:00"
$ws.Range("AA36").Value = "'2023-09-19"
$ws.Range("AA36").ClearFormats()
$ws.Range("AB36").Value = "00:00"
$ws.Range("AD36").Value = $false
$ws.Range("AE36").Value = $false
$ws.Range("AF36").Value = "'"
$ws.Range("AF36").ClearFormats()
$ws.Range("AG36").Value = $false
$ws.Range("AT36").Value = "'"
$ws.Range("AT36").ClearFormats()
$ws.Range("AW36").Value = "Larsgunnar Nilsson"
$ws.Range("AX36").Value = "Larsgunnar Nilsson, Gunilla Nilsson, Ingvor Kasselstrand, Magnus Kasselstrand"
$ws.Range("AY36").Value = "'"
$ws.Range("AY36").ClearFormats()

# Row 37
$ws.Range("A37").Value = 112204281
$ws.Range("B37").Value = 96348
$ws.Range("C37").Value = "Ovaliderad"
$ws.Range("D37").Value = "VU"
$ws.Range("E37").Value = 220787
$ws.Range("F37").Value = "Knärot"
$ws.Range("G37").Value = "Goodyera repens"
$ws.Range("H37").Value = "(L.) R. Br."
$ws.Range("I37").Value = "'10"
$ws.Range("I37").ClearFormats()
$ws.Range("J37").Value = "plantor/tuvor"
$ws.Range("K37").Value = "'"
$ws.Range("K37").ClearFormats()
$ws.Range("L37").Value = "'"
$ws.Range("L37").ClearFormats()
$ws.Range("N37").Value = "'"
$ws.Range("N37").ClearFormats()
$ws.Range("P37").Value = "Storön, Samsvik, Sm"
$ws.Range("Q37").Value = 595169.0849668512
$ws.Range("R37").Value = 6396054.017012647
$ws.Range("S37").Value = 5
$ws.Range("T37").Value = "Kalmar"
$ws.Range("U37").Value = "Västervik"
$ws.Range("V37").Value = "Småland"
$ws.Range("W37").Value = "Gladhammar"
$ws.Range("Y37").Value = "'2023-09-19"
$ws.Range("Y37").ClearFormats()
$ws.Range("Z37").Value = "00:00"
$ws.Range("AA37").Value = "'2023-09-19"
$ws.Range("AA37").ClearFormats()
$ws.Range("AB37").Value = "00:00"
$ws.Range("AD37").Value = $false
$ws.Range("AE37").Value = $false
$ws.Range("AF37").Value = "'"
$ws.Range("AF37").ClearFormats()
$ws.Range("AG37").Value = $false
$ws.Range("AT37").Value = "'"
$ws.Range("AT37").ClearFormats()
$ws.Range("AW37").Value = "Larsgunnar Nilsson"
$ws.Range("AX37").Value = "Larsgunnar Nilsson, Gunilla Nilsson, Ingvor Kasselstrand, Magnus Kasselstrand"
$ws.Range("AY37").Value = "'"
$ws.Range("AY37").ClearFormats()

# Row 38
$ws.Range("A38").Value = 112203709
$ws.Range("B38").Value = 88934
$ws.Range("C38").Value = "Ovaliderad"
$ws.Range("D38").Value = "LC"
$ws.Range("E38").Value = 5741
$ws.Range("F38").Value = "Tjockfotad fingersvamp"
$ws.Range("G38").Value = "Ramaria flavescens"
$ws.Range("H38").Value = "(Schaeff.) R. H. Petersen"
$ws.Range("I38").Value = "'4"
$ws.Range("I38").ClearFormats()
$ws.Range("J38").Value = "fruktkroppar"
$ws.Range("K38").Value = "'"
$ws.Range("K38").ClearFormats()
$ws.Range("N38").Value = "'"
$ws.Range("N38").ClearFormats()
$ws.Range("P38").Value = "A 30779, Storön, Sm"
$ws.Range("Q38").Value = 594781.2239072464
$ws.Range("R38").Value = 6396169.468659882
$ws.Range("S38").Value = 10
$ws.Range("T38").Value = "Kalmar"
$ws.Range("U38").Value = "Västervik"
$ws.Range("V38").Value = "Småland"
$ws.Range("W38").Value = "Gladhammar"
$ws.Range("Y38").Value = "'2023-09-19"
$ws.Range("Y38").ClearFormats()
$ws.Range("Z38").Value = "00:00"
$ws.Range("AA38").Value = "'2023-09-19"
$ws.Range("AA38").ClearFormats()
$ws.Range("AB38").Value = "00:00"
$ws.Range("AD38").Value = $false
$ws.Range("AE38").Value = $false
$ws.Range("AF38").Value = "'"
$ws.Range("AF38").ClearFormats()
$ws.Range("AG38").Value = $false
$ws.Range("AT38").Value = "'"
$ws.Range("AT38").ClearFormats()
$ws.Range("AW38").Value = "Magnus Kasselstrand"
$ws.Range("AX38").Value = "Magnus Kasselstrand, Ingvor Kasselstrand, Gunilla Nilsson, Larsgunnar Nilsson"
$ws.Range("AY38").Value = "'"
$ws.Range("AY38").ClearFormats()

